$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: capture the text currently sitting in the "footer" row (row 9) so
# we can re-write it one row lower (row 10) once the new product row has
# been inserted above it.
# ---------------------------------------------------------------------------
$footerDate   = $ws.Range("A9").Value2
$footerPage   = $ws.Range("G9").Value2
$footerAuthor = $ws.Range("K9").Value2

# ---------------------------------------------------------------------------
# Step 2: remove the merges that will be restructured so we can freely write
# into the cells that are about to change shape.
# ---------------------------------------------------------------------------
$ws.Range("N8:Q8").UnMerge()
$ws.Range("A9:F9").UnMerge()
$ws.Range("G9:I9").UnMerge()
$ws.Range("K9:Q9").UnMerge()

# ---------------------------------------------------------------------------
# Step 3: push the old footer row (row 9) down to row 10, copying formats
# from the originals, then write the captured values back in.
# ---------------------------------------------------------------------------
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Range("G9:I9").Copy()
$ws.Range("G10:I10").PasteSpecial(-4122)
$ws.Range("J9").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("K9:Q9").Copy()
$ws.Range("K10:Q10").PasteSpecial(-4122)

$ws.Range("A10").Value = $footerDate
$ws.Range("G10").Value = $footerPage
$ws.Range("K10").Value = $footerAuthor

$ws.Rows(10).RowHeight = 16.5
$ws.Range("A10:F10").Merge()
$ws.Range("G10:I10").Merge()
$ws.Range("K10:Q10").Merge()

# ---------------------------------------------------------------------------
# Step 4: push the "totals" row (old row 8) down to row 9, keeping its
# format, and give it the new total value (136).
# ---------------------------------------------------------------------------
$ws.Range("N8:Q8").Copy()
$ws.Range("N9:Q9").PasteSpecial(-4122)
$ws.Range("N9").Value = 136

$ws.Rows(9).RowHeight = 26.25
$ws.Range("N9:Q9").Merge()

# clear the leftover values that used to live in row 8 before we reuse it
$ws.Range("N8:Q8").ClearContents()

# ---------------------------------------------------------------------------
# Step 5: build the new data row 8 (second sale line), copying the cell
# formatting from row 7 (the first sale line) column by column.
# ---------------------------------------------------------------------------
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$ws.Range("C7:G7").Copy()
$ws.Range("C8:G8").PasteSpecial(-4122)
$ws.Range("H7:K7").Copy()
$ws.Range("H8:K8").PasteSpecial(-4122)
$ws.Range("L7:M7").Copy()
$ws.Range("L8:M8").PasteSpecial(-4122)
$ws.Range("N7:O7").Copy()
$ws.Range("N8:O8").PasteSpecial(-4122)
$ws.Range("P7").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("Q7").Copy()
$ws.Range("Q8").PasteSpecial(-4122)

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "TOP ZED RESMOOTH GEL"
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = "0"
$ws.Range("N8").Value = "49.00"
$ws.Range("P8").Value = "49.0000"
$ws.Range("Q8").Value = "1:0"

$ws.Rows(8).RowHeight = 24.75
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# ---------------------------------------------------------------------------
# Step 6: fill in the first data row (row 7) - formats already in place.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "FLECTOR 50MG 30 CAPS"
$ws.Range("H7").Value = "0:2"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "87.00"
$ws.Range("P7").Value = "87.0000"
$ws.Range("Q7").Value = "1:0"
